$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.309.02'
$ws.Range('E2').Value = '  +3.01%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.005.75'
$ws.Range('E3').Value = '  +7.13%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7893'
$ws.Range('E5').Value = '  +67.83%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '259.57'
$ws.Range('E6').Value = '  +6.19%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3610'
$ws.Range('E8').Value = '  +25.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '28.46'
$ws.Range('E9').Value = '  +29.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07079'
$ws.Range('E10').Value = '  +9.03%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8607'
$ws.Range('E11').Value = '  +17.84%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08191'

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.009.15'
$ws.Range('E13').Value = '  +7.34%  '

$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '101.55'
$ws.Range('E14').Value = '  +1.64%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.635'
$ws.Range('E15').Value = '  +8.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '276.65'
$ws.Range('E16').Value = '  -3.75%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.26'
$ws.Range('E17').Value = '  +16.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '31.322.55'
$ws.Range('E18').Value = '  +3.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.955'
$ws.Range('E19').Value = '  +12.36%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000008017'
$ws.Range('E20').Value = '  +7.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.271.39'
$ws.Range('E21').Value = '  +7.54%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.0000'
$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.202'
$ws.Range('E24').Value = '  +13.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.15'
$ws.Range('E25').Value = '  +12.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1518'
$ws.Range('E26').Value = '  +57.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.17'
$ws.Range('E27').Value = '  +1.72%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.11'
$ws.Range('E28').Value = '  +5.87%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.403'
$ws.Range('E29').Value = '  +26.76%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.625'
$ws.Range('E30').Value = '  +9.40%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.632'
$ws.Range('E31').Value = '  +9.39%  '

$ws.Range('E32').Value = '  +3.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.448'
$ws.Range('E33').Value = '  +7.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05243'
$ws.Range('E34').Value = '  +8.88%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7809'
$ws.Range('E35').Value = '  +13.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.223'
$ws.Range('E36').Value = '  +8.59%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.812'
$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02013'
$ws.Range('E38').Value = '  +5.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.938'
$ws.Range('E39').Value = '  +3.36%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.731'
$ws.Range('E40').Value = '  +7.26%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.51'
$ws.Range('E41').Value = '  +5.74%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4776'
$ws.Range('E42').Value = '  +12.85%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.166'
$ws.Range('E43').Value = '  +10.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '108.19'
$ws.Range('E44').Value = '  +7.08%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8586'
$ws.Range('E45').Value = '  +4.23%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.861'
$ws.Range('E46').Value = '  +12.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.965'
$ws.Range('E48').Value = '  +2.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4378'
$ws.Range('E49').Value = '  +12.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.98'
$ws.Range('E50').Value = '  +5.60%  '

$ws.Range('E51').Value = '  +14.89%  '
